# Update Quiz Dashboard scripts
# - Fix method name casing: verifyQuizDashBoardClasses -> verifyQuizDashboardClasses
#   for all "Quiz Dashboard Classes" rows (11-19) in column C.
# - Fill in column F (Value) for rows 17-19 with "Class 6-C, Class 7-A".
# - Update the active selection to F17:F19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the method name casing in column C for rows 11 through 19.
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 3).Value = "verifyQuizDashboardClasses"
}

# Populate column F for rows 17-19 (previously blank).
$ws.Range("F17").Value = "Class 6-C, Class 7-A"
$ws.Range("F18").Value = "Class 6-C, Class 7-A"
$ws.Range("F19").Value = "Class 6-C, Class 7-A"

# Update the selected range shown in the sheet view.
$ws.Range("F17:F19").Select()
